$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Minor numeric revisions in column C (and B84) for existing rows ---
$ws.Range("C39").Value = 87.48
$ws.Range("C40").Value = 88.36
$ws.Range("C41").Value = 89.31999999999999
$ws.Range("C44").Value = 92.23
$ws.Range("C45").Value = 91.5
$ws.Range("C46").Value = 93
$ws.Range("C49").Value = 95.34
$ws.Range("C50").Value = 95.23999999999999
$ws.Range("C54").Value = 100.61
$ws.Range("C59").Value = 101.56
$ws.Range("C62").Value = 101.3
$ws.Range("C66").Value = 100.79
$ws.Range("C68").Value = 98.38
$ws.Range("C70").Value = 94.73999999999999
$ws.Range("C71").Value = 97.66
$ws.Range("C72").Value = 100.94
$ws.Range("C74").Value = 102.49
$ws.Range("C76").Value = 101.29
$ws.Range("C77").Value = 102.78
$ws.Range("C78").Value = 100.16
$ws.Range("C79").Value = 102.17
$ws.Range("C80").Value = 102.02
$ws.Range("C82").Value = 101.81
$ws.Range("C83").Value = 97.84999999999999
$ws.Range("B84").Value = 98.41
$ws.Range("C84").Value = 99.02
$ws.Range("C85").Value = 101.08
$ws.Range("C86").Value = 102.13
$ws.Range("C87").Value = 102.93

# --- Append the new quarterly data row (row 88) ---
# Use a quote-prefixed text entry so Excel stores the date label as a
# plain string (matching the existing A-column pattern) instead of
# auto-converting it to a date serial, then clear the resulting
# quote-prefix style so the cell carries no style override.
$ws.Range("A88").Value = "'01-07-2021"
$ws.Range("A88").Style = "Normal"
$ws.Range("B88").Value = 101.29
$ws.Range("C88").Value = 101.87
